$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 581
"done"
